$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update "Story points" (col B) and "Total Points" (col C) for rows 2-14 ---
$ws.Range("B2").Value = 60
$ws.Range("C2").Value = 50

$ws.Range("B3").Value = 55
$ws.Range("C3").Value = 50

$ws.Range("B4").Value = 50
$ws.Range("C4").Value = 50

$ws.Range("B5").Value = 45
$ws.Range("C5").Value = 50

$ws.Range("B6").Value = 40
$ws.Range("C6").Value = 50

$ws.Range("B7").Value = 35
$ws.Range("C7").Value = 50

$ws.Range("B8").Value = 30
$ws.Range("C8").Value = 50

$ws.Range("B9").Value = 25
$ws.Range("C9").Value = 50

$ws.Range("B10").Value = 20
$ws.Range("C10").Value = 50

# Rows 11-14 also change "Completed" (col D)
$ws.Range("B11").Value = 15
$ws.Range("C11").Value = 50
$ws.Range("D11").Value = 47

$ws.Range("B12").Value = 10
$ws.Range("C12").Value = 50
$ws.Range("D12").Value = 47

$ws.Range("B13").Value = 5
$ws.Range("C13").Value = 50
$ws.Range("D13").Value = 47

$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 50
$ws.Range("D14").Value = 50

# --- Rows 15 and 16 are no longer part of the data range: clear them out ---
$ws.Range("A15:E16").ClearContents()
